# Auto-generated edit script applying the diff cell-by-cell.
# Re-creates each changed cell as a text (inline-string-equivalent) value,
# using a leading apostrophe to force text interpretation (avoiding numeric/percent
# auto-conversion), then resets the style back to "Normal" so no stray number format
# style gets attached to the cell (matching the original un-styled cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.53%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.13%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.259"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.92%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05704"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.08%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.05%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.191"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.98%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.79%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.99%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'One"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.01005"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1,578.09%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'WazirX"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1371"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.68%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.07041"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.52%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.03173"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'8.63%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09282"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.10%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.001528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.59%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.006074"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.99%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.493"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.09%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'BTSEToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'2.230"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.65%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3160"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.46%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.03293"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.00%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-3.65%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.487"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.24%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04090"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.04%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.07%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.23%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004143"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-17.66%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.86%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001449"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-25.27%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03756"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.17%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1066"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.41%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.003719"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-35.44%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002448"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'22.37%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009374"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-5.93%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005247"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.70%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.07499"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'24.96%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002444"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-4.74%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E50").Style = "Normal"
